$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting the existing
# Model / Passenger Capacity / Full Tanks Full data to columns B:D
$ws.Columns.Item(1).Insert()

# Copy the bold/filled header style from the shifted header (now B1) onto
# the new "NO AIRCRAFT" header cell in A1
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New header label
$ws.Range("A1").Value = "NO AIRCRAFT"

# Sequential aircraft numbers for the data rows
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6

# Re-normalize the last row's numeric values (carried over as "600.0")
$ws.Range("C7").Value = 0
$ws.Range("C7").Value = 600
$ws.Range("D7").Value = 0
$ws.Range("D7").Value = 600

# Column width adjustments made while laying out the extra
# airport/aircraft selection area to the right of the table
$ws.Columns.Item(2).ColumnWidth = 34.917
$ws.Columns.Item(3).ColumnWidth = 13.584
$ws.Columns.Item(4).ColumnWidth = 14.751
$ws.Columns.Item(5).ColumnWidth = 13.584
$ws.Columns.Item(7).ColumnWidth = 16.417
$ws.Columns.Item(8).ColumnWidth = 13.417
$ws.Columns.Item(10).ColumnWidth = 18.917

# Restore the selection used while building the new layout
$ws.Range("E1:H7").Select()

Write-Host "done"
